# Applies the numeric cell updates from the scheduled runner diff
# across the ALC/ARM/CRP/CUL/GSM/LTW/WVR sheets of the Excalibur_Profits workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4566.125
$ws.Range("I18").Value = 755
$ws.Range("K18").Value = 755
$ws.Range("M18").Value = -471
$ws.Range("H46").Value = 524.8333
$ws.Range("I46").Value = 529.8
$ws.Range("K46").Value = 1589.4
$ws.Range("M46").Value = -1470.4
$ws.Range("H53").Value = 372
$ws.Range("I53").Value = 358
$ws.Range("J53").Value = 382.5
$ws.Range("K53").Value = 358
$ws.Range("L53").Value = 382.5
$ws.Range("M53").Value = 279
$ws.Range("N53").Value = -1656.5
$ws.Range("H60").Value = 524.8333
$ws.Range("I60").Value = 529.8
$ws.Range("K60").Value = 1589.4
$ws.Range("M60").Value = -1105.4
$ws.Range("H76").Value = 6912.1665
$ws.Range("I76").Value = 4983
$ws.Range("J76").Value = 7555.222
$ws.Range("K76").Value = 4983
$ws.Range("L76").Value = 7555.222
$ws.Range("M76").Value = -4668
$ws.Range("N76").Value = -8185.222
$ws.Range("H79").Value = 6912.1665
$ws.Range("I79").Value = 4983
$ws.Range("J79").Value = 7555.222
$ws.Range("K79").Value = 4983
$ws.Range("L79").Value = 7555.222
$ws.Range("M79").Value = -3891
$ws.Range("N79").Value = -9739.222
$ws.Range("H86").Value = 3586.8857
$ws.Range("I86").Value = 2208.8333
$ws.Range("J86").Value = 5046
$ws.Range("K86").Value = 2208.8333
$ws.Range("L86").Value = 5046
$ws.Range("M86").Value = -1085.8333
$ws.Range("N86").Value = -7292
$ws.Range("H89").Value = 3586.8857
$ws.Range("I89").Value = 2208.8333
$ws.Range("J89").Value = 5046
$ws.Range("K89").Value = 11044.1665
$ws.Range("L89").Value = 25230
$ws.Range("M89").Value = -5428.166499999999
$ws.Range("N89").Value = -36462
$ws.Range("H94").Value = 6413.923
$ws.Range("I94").Value = 2781.75
$ws.Range("J94").Value = 50000
$ws.Range("K94").Value = 2781.75
$ws.Range("L94").Value = 50000
$ws.Range("M94").Value = -2330.75
$ws.Range("N94").Value = -50902

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 2210
$ws.Range("I12").Value = 1716.6666
$ws.Range("J12").Value = 2950
$ws.Range("K12").Value = 1716.6666
$ws.Range("L12").Value = 2950
$ws.Range("M12").Value = -1543.6666
$ws.Range("N12").Value = -3296
$ws.Range("H31").Value = 7063.3335
$ws.Range("I31").Value = 7063.3335
$ws.Range("K31").Value = 7063.3335
$ws.Range("M31").Value = -6769.3335
$ws.Range("H45").Value = 3278.95
$ws.Range("I45").Value = 3024.4285
$ws.Range("K45").Value = 3024.4285
$ws.Range("M45").Value = -2647.4285
$ws.Range("H102").Value = 57355.5
$ws.Range("I102").Value = 57355.5
$ws.Range("K102").Value = 57355.5
$ws.Range("M102").Value = -55733.5
$ws.Range("H110").Value = 1983
$ws.Range("I110").Value = 950
$ws.Range("J110").Value = 2499.5
$ws.Range("K110").Value = 950
$ws.Range("L110").Value = 2499.5
$ws.Range("M110").Value = 1095
$ws.Range("N110").Value = -6589.5
$ws.Range("H122").Value = 3543
$ws.Range("I122").Value = 2857.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 8572.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -6122.5
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1055883.2
$ws.Range("I132").Value = 1430834.4
$ws.Range("K132").Value = 4292503.199999999
$ws.Range("M132").Value = -4289973.199999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 655.75
$ws.Range("I16").Value = 463.7143
$ws.Range("K16").Value = 463.7143
$ws.Range("M16").Value = -176.7143
$ws.Range("H107").Value = 539.79486
$ws.Range("I107").Value = 488.33334
$ws.Range("K107").Value = 488.33334
$ws.Range("M107").Value = 1431.66666
$ws.Range("H113").Value = 655.75
$ws.Range("I113").Value = 463.7143
$ws.Range("K113").Value = 463.7143
$ws.Range("M113").Value = 1706.2857
$ws.Range("H134").Value = 1794.2
$ws.Range("I134").Value = 1465
$ws.Range("J134").Value = 2699.5
$ws.Range("K134").Value = 4395
$ws.Range("L134").Value = 8098.5
$ws.Range("M134").Value = -1860
$ws.Range("N134").Value = -13168.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1244.8334
$ws.Range("J98").Value = 1242.5
$ws.Range("L98").Value = 3727.5
$ws.Range("N98").Value = -6723.5
$ws.Range("H102").Value = 8553.333000000001
$ws.Range("J102").Value = 8553.333000000001
$ws.Range("L102").Value = 25659.999
$ws.Range("N102").Value = -30527.999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 41432.918
$ws.Range("J92").Value = 41432.918
$ws.Range("L92").Value = 41432.918
$ws.Range("N92").Value = -45176.918
$ws.Range("H118").Value = 31578
$ws.Range("J118").Value = 31578
$ws.Range("L118").Value = 31578
$ws.Range("N118").Value = -34892
$ws.Range("H122").Value = 66504.5
$ws.Range("I122").Value = 70605.47
$ws.Range("K122").Value = 211816.41
$ws.Range("M122").Value = -209366.41

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 116821.22
$ws.Range("I20").Value = 49999
$ws.Range("J20").Value = 125174
$ws.Range("K20").Value = 49999
$ws.Range("L20").Value = 125174
$ws.Range("M20").Value = -49773
$ws.Range("N20").Value = -125626
$ws.Range("H40").Value = 2296.96
$ws.Range("I40").Value = 2331.261
$ws.Range("J40").Value = 1902.5
$ws.Range("K40").Value = 2331.261
$ws.Range("L40").Value = 1902.5
$ws.Range("M40").Value = -2195.261
$ws.Range("N40").Value = -2174.5
$ws.Range("H68").Value = 1910.4445
$ws.Range("I68").Value = 1827.7142
$ws.Range("K68").Value = 1827.7142
$ws.Range("M68").Value = -1078.7142
$ws.Range("H71").Value = 1910.4445
$ws.Range("I71").Value = 1827.7142
$ws.Range("K71").Value = 9138.571
$ws.Range("M71").Value = -5394.571
$ws.Range("H93").Value = 974.1429000000001
$ws.Range("I93").Value = 797.6667
$ws.Range("K93").Value = 797.6667
$ws.Range("M93").Value = 450.3333
$ws.Range("H132").Value = 16292.579
$ws.Range("I132").Value = 21404.615
$ws.Range("K132").Value = 64213.845
$ws.Range("M132").Value = -61683.845
$ws.Range("H136").Value = 2988.8
$ws.Range("I136").Value = 2280.5
$ws.Range("K136").Value = 6841.5
$ws.Range("M136").Value = -4291.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1718.4517
$ws.Range("J113").Value = 3895.7273
$ws.Range("L113").Value = 11687.1819
$ws.Range("N113").Value = -16027.1819
$ws.Range("H132").Value = 2383015
$ws.Range("I132").Value = 2925834.2
$ws.Range("K132").Value = 8777502.600000001
$ws.Range("M132").Value = -8774972.600000001
$ws.Range("H136").Value = 9153.23
$ws.Range("I136").Value = 9039.360000000001
$ws.Range("K136").Value = 27118.08
$ws.Range("M136").Value = -24568.08

Write-Output "Updated 176 cells across 7 sheets"
